$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 1
$ws.Range("I8").Value = 1
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 136
$ws.Range("H19").Value = 2630.9167
$ws.Range("J19").Value = 2630.9167
$ws.Range("L19").Value = 2630.9167
$ws.Range("N19").Value = -2980.9167
$ws.Range("H32").Value = 4944.1113
$ws.Range("I32").Value = 1357
$ws.Range("J32").Value = 17499
$ws.Range("K32").Value = 1357
$ws.Range("L32").Value = 17499
$ws.Range("M32").Value = -1031
$ws.Range("N32").Value = -18151
$ws.Range("H33").Value = 2662.2222
$ws.Range("I33").Value = 3308.9443
$ws.Range("J33").Value = 1368.7778
$ws.Range("K33").Value = 3308.9443
$ws.Range("L33").Value = 1368.7778
$ws.Range("M33").Value = -3079.9443
$ws.Range("N33").Value = -1826.7778
$ws.Range("H38").Value = 1890.8572
$ws.Range("J38").Value = 2907.8333
$ws.Range("L38").Value = 8723.499899999999
$ws.Range("N38").Value = -9467.499899999999
$ws.Range("H39").Value = 3513.4167
$ws.Range("I39").Value = 560
$ws.Range("J39").Value = 4990.125
$ws.Range("K39").Value = 1680
$ws.Range("L39").Value = 14970.375
$ws.Range("M39").Value = -1384
$ws.Range("N39").Value = -15562.375
$ws.Range("H45").Value = 6999.25
$ws.Range("J45").Value = 9999
$ws.Range("L45").Value = 29997
$ws.Range("N45").Value = -30381
$ws.Range("H74").Value = 13834.529
$ws.Range("I74").Value = 13834.529
$ws.Range("K74").Value = 13834.529
$ws.Range("M74").Value = -12898.529
$ws.Range("H76").Value = 14716.5
$ws.Range("I76").Value = 19112.125
$ws.Range("K76").Value = 19112.125
$ws.Range("M76").Value = -18797.125
$ws.Range("H77").Value = 13834.529
$ws.Range("I77").Value = 13834.529
$ws.Range("K77").Value = 69172.645
$ws.Range("M77").Value = -64492.645
$ws.Range("H79").Value = 14716.5
$ws.Range("I79").Value = 19112.125
$ws.Range("K79").Value = 19112.125
$ws.Range("M79").Value = -18020.125
$ws.Range("H86").Value = 1628.875
$ws.Range("I86").Value = 1752.1818
$ws.Range("J86").Value = 1357.6
$ws.Range("K86").Value = 1752.1818
$ws.Range("L86").Value = 1357.6
$ws.Range("M86").Value = -629.1818000000001
$ws.Range("N86").Value = -3603.6
$ws.Range("H89").Value = 1628.875
$ws.Range("I89").Value = 1752.1818
$ws.Range("J89").Value = 1357.6
$ws.Range("K89").Value = 8760.909
$ws.Range("L89").Value = 6788
$ws.Range("M89").Value = -3144.909
$ws.Range("N89").Value = -18020
$ws.Range("H137").Value = 1349.9474
$ws.Range("I137").Value = 1324.9445
$ws.Range("K137").Value = 3974.8335
$ws.Range("M137").Value = -1424.8335

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 5858.1
$ws.Range("I63").Value = 1716.3
$ws.Range("J63").Value = 9999.9
$ws.Range("K63").Value = 1716.3
$ws.Range("L63").Value = 9999.9
$ws.Range("M63").Value = -1030.3
$ws.Range("N63").Value = -11371.9
$ws.Range("H66").Value = 5858.1
$ws.Range("I66").Value = 1716.3
$ws.Range("J66").Value = 9999.9
$ws.Range("K66").Value = 8581.5
$ws.Range("L66").Value = 49999.5
$ws.Range("M66").Value = -5149.5
$ws.Range("N66").Value = -56863.5
$ws.Range("H74").Value = 4043.7693
$ws.Range("I74").Value = 3339.25
$ws.Range("K74").Value = 3339.25
$ws.Range("M74").Value = -2465.25
$ws.Range("H77").Value = 4043.7693
$ws.Range("I77").Value = 3339.25
$ws.Range("K77").Value = 16696.25
$ws.Range("M77").Value = -12328.25
$ws.Range("H97").Value = 674.75
$ws.Range("I97").Value = 696.92
$ws.Range("J97").Value = 490
$ws.Range("K97").Value = 696.92
$ws.Range("L97").Value = 490
$ws.Range("M97").Value = -200.92
$ws.Range("N97").Value = -1482

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 438.16666
$ws.Range("J64").Value = 515.5
$ws.Range("L64").Value = 515.5
$ws.Range("N64").Value = -965.5
$ws.Range("H67").Value = 438.16666
$ws.Range("J67").Value = 515.5
$ws.Range("L67").Value = 515.5
$ws.Range("N67").Value = -2075.5
$ws.Range("H86").Value = 7577487.5
$ws.Range("I86").Value = 10418340
$ws.Range("K86").Value = 10418340
$ws.Range("M86").Value = -10417217
$ws.Range("H89").Value = 7577487.5
$ws.Range("I89").Value = 10418340
$ws.Range("K89").Value = 52091700
$ws.Range("M89").Value = -52086084
$ws.Range("H94").Value = 7578644.5
$ws.Range("I94").Value = 14707148
$ws.Range("K94").Value = 14707148
$ws.Range("M94").Value = -14706697
$ws.Range("H105").Value = 1850
$ws.Range("I105").Value = 1850
$ws.Range("K105").Value = 1850
$ws.Range("M105").Value = -103
$ws.Range("H134").Value = 3458.9583
$ws.Range("I134").Value = 3391.9565
$ws.Range("K134").Value = 10175.8695
$ws.Range("M134").Value = -7640.869499999999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 395.5
$ws.Range("I5").Value = 277.18182
$ws.Range("K5").Value = 831.54546
$ws.Range("M5").Value = -719.54546
$ws.Range("H38").Value = 123.5
$ws.Range("J38").Value = 146.5
$ws.Range("L38").Value = 439.5
$ws.Range("N38").Value = -1133.5
$ws.Range("H56").Value = 3960.3572
$ws.Range("I56").Value = 3960.3572
$ws.Range("K56").Value = 3960.3572
$ws.Range("M56").Value = -3430.3572
$ws.Range("H75").Value = 497.5
$ws.Range("J75").Value = 497.5
$ws.Range("L75").Value = 1492.5
$ws.Range("N75").Value = -3488.5
$ws.Range("H78").Value = 497.5
$ws.Range("J78").Value = 497.5
$ws.Range("L78").Value = 4477.5
$ws.Range("N78").Value = -14461.5
$ws.Range("H113").Value = 1254.3158
$ws.Range("I113").Value = 573.9286
$ws.Range("K113").Value = 1721.7858
$ws.Range("M113").Value = 448.2142000000001
$ws.Range("H120").Value = 15995.6
$ws.Range("I120").Value = 8326.666999999999
$ws.Range("K120").Value = 24980.001
$ws.Range("M120").Value = -20142.001
$ws.Range("H130").Value = 5999
$ws.Range("J130").Value = 5999
$ws.Range("L130").Value = 17997
$ws.Range("N130").Value = -28037
$ws.Range("H135").Value = 395.5
$ws.Range("I135").Value = 277.18182
$ws.Range("K135").Value = 2494.63638
$ws.Range("M135").Value = 40.36362000000008

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6667.25
$ws.Range("I70").Value = 6667.25
$ws.Range("K70").Value = 6667.25
$ws.Range("M70").Value = -6397.25
$ws.Range("H73").Value = 6667.25
$ws.Range("I73").Value = 6667.25
$ws.Range("K73").Value = 6667.25
$ws.Range("M73").Value = -5731.25
$ws.Range("H97").Value = 226
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()
$ws.Range("H102").Value = 2522.2104
$ws.Range("I102").Value = 2557.2307
$ws.Range("K102").Value = 2557.2307
$ws.Range("M102").Value = -935.2307000000001

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1900.1666
$ws.Range("I82").Value = 1998.2858
$ws.Range("J82").Value = 1837.7273
$ws.Range("K82").Value = 1998.2858
$ws.Range("L82").Value = 1837.7273
$ws.Range("M82").Value = -1637.2858
$ws.Range("N82").Value = -2559.7273
$ws.Range("H85").Value = 1900.1666
$ws.Range("I85").Value = 1998.2858
$ws.Range("J85").Value = 1837.7273
$ws.Range("K85").Value = 1998.2858
$ws.Range("L85").Value = 1837.7273
$ws.Range("M85").Value = -750.2858000000001
$ws.Range("N85").Value = -4333.7273
$ws.Range("H132").Value = 10979.777
$ws.Range("I132").Value = 8743.294
$ws.Range("K132").Value = 26229.882
$ws.Range("M132").Value = -23699.882

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 68799.5
$ws.Range("J74").Value = 68799.5
$ws.Range("L74").Value = 68799.5
$ws.Range("N74").Value = -70671.5
$ws.Range("H77").Value = 68799.5
$ws.Range("J77").Value = 68799.5
$ws.Range("L77").Value = 206398.5
$ws.Range("N77").Value = -215758.5
$ws.Range("H126").Value = 3173.8823
$ws.Range("I126").Value = 2950.2666
$ws.Range("J126").Value = 4851
$ws.Range("K126").Value = 8850.799800000001
$ws.Range("L126").Value = 14553
$ws.Range("M126").Value = -6380.799800000001
$ws.Range("N126").Value = -19493
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()
$ws.Range("H136").Value = 3688.5925
$ws.Range("I136").Value = 3483.68
$ws.Range("K136").Value = 10451.04
$ws.Range("M136").Value = -7901.039999999999
